# "keeping track of mpcs" - add newly tracked MPC variables/quarters to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Fill-Row($row, $startCol, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $startCol + $i).Value = $values[$i]
    }
}

# ---------------------------------------------------------------------------
# Row 9 (federal_ui_arp): fill in the per-quarter MPC decay values
# ---------------------------------------------------------------------------
Fill-Row 9 5 @(0.2, 0.17, 0.16, 0.15, 0.09, 0.05, 0.05, 0.04)

# ---------------------------------------------------------------------------
# Row 10: state_ui_arp
# ---------------------------------------------------------------------------
$ws.Cells.Item(10, 1).Value = "state_ui_arp"
$ws.Cells.Item(10, 2).Value = "mpc_vulnerable_arp(.x)"
$ws.Cells.Item(10, 3).Value = "246-254"
$ws.Cells.Item(10, 4).Value = "code"
Fill-Row 10 5 @(0.2, 0.17, 0.16, 0.15, 0.09, 0.05, 0.05, 0.04)

# ---------------------------------------------------------------------------
# Row 11: federal_other_vulnerable_arp
# ---------------------------------------------------------------------------
$ws.Cells.Item(11, 1).Value = "federal_other_vulnerable_arp"
$ws.Cells.Item(11, 2).Value = "mpc_vulnerable_arp(.x)"
$ws.Cells.Item(11, 3).Value = "246-255"
$ws.Cells.Item(11, 4).Value = "code"
Fill-Row 11 5 @(0.2, 0.17, 0.16, 0.15, 0.09, 0.05, 0.05, 0.04)

# ---------------------------------------------------------------------------
# Row 12: rebate_checks_arp
# ---------------------------------------------------------------------------
$directAidCode = "    across(`n      .cols = all_of(`n        c(""rebate_checks_arp"", ""federal_other_direct_aid_arp"") %>% paste0(""_minus_neutral"")`n      ),`n      #same as above, applying a different MPC function to these `n      .fns = ~ mpc_direct_aid_arp(.),`n      .names = ""{.col}_post_mpc""`n    ),"

$ws.Cells.Item(12, 1).Value = "rebate_checks_arp"
$ws.Cells.Item(12, 2).Value = "mpc_direct_aid_arp(.)"
$ws.Cells.Item(12, 3).Value = "254-261"
$ws.Cells.Item(12, 4).Value = $directAidCode
Fill-Row 12 5 @(0.14, 0.1, 0.1, 0.05, 0.05, 0.05, 0.05, 0.05, 0.03, 0.03, 0.03, 0.03)

# ---------------------------------------------------------------------------
# Row 13: federal_other_direct_aid_arp
# ---------------------------------------------------------------------------
$ws.Cells.Item(13, 1).Value = "federal_other_direct_aid_arp"
$ws.Cells.Item(13, 2).Value = "mpc_direct_aid_arp(.)"
$ws.Cells.Item(13, 3).Value = "254-261"
$ws.Cells.Item(13, 4).Value = $directAidCode
Fill-Row 13 5 @(0.14, 0.1, 0.1, 0.05, 0.05, 0.05, 0.05, 0.05, 0.05, 0.03, 0.03, 0.03)

# ---------------------------------------------------------------------------
# Row 14: federal_student_loans
# ---------------------------------------------------------------------------
$studentLoansCode = "    across(`n      .cols = all_of(`n        c(""federal_student_loans"") %>% paste0(""_minus_neutral"")`n      ),`n      #same as above, applying a different MPC function to these `n      .fns = ~ mpc_student_loans(.),`n      .names = ""{.col}_post_mpc""`n    ),"

$ws.Cells.Item(14, 1).Value = "federal_student_loans"
$ws.Cells.Item(14, 2).Value = "mpc_student_loans(.)"
$ws.Cells.Item(14, 3).Value = "262-269"
$ws.Cells.Item(14, 4).Value = $studentLoansCode
Fill-Row 14 5 @(0.2, 0.17, 0.16, 0.15, 0.09, 0.05, 0.05, 0.04)

# ---------------------------------------------------------------------------
# Row 15: supply_side_ira
# ---------------------------------------------------------------------------
$supplySideIraCode = "    across(`n      .cols = any_of(`n        c(""supply_side_ira"") %>% paste0(""_minus_neutral"")`n      ),`n      #getting the post mpc levels for the ARP variables`n      .fns = ~ mpc_supply_side_ira(.x),`n      .names = ""{.col}_post_mpc""`n    ),"

$ws.Cells.Item(15, 1).Value = "supply_side_ira"
$ws.Cells.Item(15, 2).Value = "mpc_supply_side_ira(.x)"
$ws.Cells.Item(15, 3).Value = "270-277"
$ws.Cells.Item(15, 4).Value = $supplySideIraCode
$ws.Cells.Item(15, 5).Value = 1

# ---------------------------------------------------------------------------
# Row 16: federal_aid_to_small_businesses
# ---------------------------------------------------------------------------
$smallBusinessCode = "    #same as above, applying a different MPC function to this`n    federal_aid_to_small_businesses_arp_minus_neutral_post_mpc = `n      mpc_small_businesses_arp ((federal_aid_to_small_businesses_arp_minus_neutral))`n  )"

$ws.Cells.Item(16, 1).Value = "federal_aid_to_small_businesses"
$ws.Cells.Item(16, 2).Value = "mpc_small_businesses_arp()"
$ws.Cells.Item(16, 3).Value = "278-281"
$ws.Cells.Item(16, 4).Value = $smallBusinessCode
Fill-Row 16 5 @(0.04, 0.04, 0.017, 0.017, 0.017, 0.017, 0.017, 0.017, 0.017, 0.017, 0.017, 0.017)

# ---------------------------------------------------------------------------
# Leave the cursor where the author left it when saving
# ---------------------------------------------------------------------------
$ws.Range("H29").Select()
